$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue 'D2' '68.973.54'
Set-TextValue 'E2' '  +3.81%  '

Set-TextValue 'D3' '3.716.74'
Set-TextValue 'E3' '  +2.57%  '

Set-TextValue 'E4' '  -0.37%  '

Set-TextValue 'D5' '612.68'
Set-TextValue 'E5' '  +10.43%  '

Set-TextValue 'D6' '193.12'
Set-TextValue 'E6' '  +16.03%  '

Set-TextValue 'D7' '0.638'
Set-TextValue 'E7' '  +5.54%  '

Set-TextValue 'E8' '  -0.86%  '

Set-TextValue 'D9' '0.725'
Set-TextValue 'E9' '  +6.55%  '

Set-TextValue 'D10' '0.161'
Set-TextValue 'E10' '  +3.88%  '

Set-TextValue 'D11' '60.05'
Set-TextValue 'E11' '  +21.55%  '

Set-TextValue 'E12' '  +4.11%  '

Set-TextValue 'D13' '10.49'
Set-TextValue 'E13' '  +3.50%  '

Set-TextValue 'D14' '4.313.07'
Set-TextValue 'E14' '  +1.42%  '

Set-TextValue 'D15' '3.719.86'
Set-TextValue 'E15' '  +1.54%  '

Set-TextValue 'B16' 'Chainlink'
Set-TextValue 'C16' 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-TextValue 'D16' '19.50'
Set-TextValue 'E16' '  +4.04%  '

Set-TextValue 'B17' 'Polygon'
Set-TextValue 'C17' 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
Set-TextValue 'D17' '1.15'
Set-TextValue 'E17' '  +6.18%  '

Set-TextValue 'E18' '  +1.29%  '

Set-TextValue 'D19' '12.96'
Set-TextValue 'E19' '  +4.75%  '

Set-TextValue 'D20' '68.827.98'
Set-TextValue 'E20' '  +3.66%  '

Set-TextValue 'D21' '412.41'
Set-TextValue 'E21' '  +5.62%  '

Set-TextValue 'D22' '4.57'
Set-TextValue 'E22' '  +5.83%  '

Set-TextValue 'D23' '89.98'
Set-TextValue 'E23' '  +6.10%  '

Set-TextValue 'D24' '3.09'
Set-TextValue 'E24' '  +5.94%  '

Set-TextValue 'B25' 'InternetComputer(DFINITY)'
Set-TextValue 'C25' 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextValue 'D25' '13.04'
Set-TextValue 'E25' '  +7.00%  '

Set-TextValue 'B26' 'RenderToken'
Set-TextValue 'C26' 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue 'D26' '11.26'
Set-TextValue 'E26' '  +11.09%  '

Set-TextValue 'D27' '3.81'
Set-TextValue 'E27' '  +6.89%  '

Set-TextValue 'D28' '6.04'
Set-TextValue 'E28' '  +1.42%  '

Set-TextValue 'D29' '9.76'
Set-TextValue 'E29' '  +7.57%  '

Set-TextValue 'D30' '32.85'
Set-TextValue 'E30' '  +4.50%  '

Set-TextValue 'D31' '7.85'
Set-TextValue 'E31' '  +8.02%  '

Set-TextValue 'D32' '12.78'
Set-TextValue 'E32' '  +6.13%  '

Set-TextValue 'D33' '646.38'
Set-TextValue 'E33' '  +14.39%  '

Set-TextValue 'E34' '  +9.38%  '

Set-TextValue 'D35' '46.05'
Set-TextValue 'E35' '  +11.65%  '

Set-TextValue 'D36' '67.56'
Set-TextValue 'E36' '  +6.96%  '

Set-TextValue 'D37' '0.0₃0836'
Set-TextValue 'E37' '  -1.73%  '

Set-TextValue 'D38' '0.417'
Set-TextValue 'E38' '  +10.11%  '

Set-TextValue 'D39' '1.00'
Set-TextValue 'E39' '  +0.29%  '

Set-TextValue 'E40' '  -0.36%  '

Set-TextValue 'E41' '  +8.89%  '

Set-TextValue 'D42' '3.07'
Set-TextValue 'E42' '  +6.69%  '

Set-TextValue 'D43' '0.0448'
Set-TextValue 'E43' '  +6.26%  '

Set-TextValue 'E44' '  +8.30%  '

Set-TextValue 'D45' '2.896.76'
Set-TextValue 'E45' '  +10.00%  '

Set-TextValue 'E46' '  +7.73%  '

Set-TextValue 'D47' '9.30'
Set-TextValue 'E47' '  +5.16%  '

Set-TextValue 'D48' '2.73'
Set-TextValue 'E48' '  +3.72%  '

Set-TextValue 'D49' '145.78'
Set-TextValue 'E49' '  +3.61%  '

Set-TextValue 'D50' '3.13'
Set-TextValue 'E50' '  +3.84%  '

Set-TextValue 'D51' '2.57'
Set-TextValue 'E51' '  -8.94%  '
